# The document ends with an empty paragraph (just before the final
# sectPr). Replace that paragraph (including its paragraph mark) with
# one that has two tab runs and a paragraph-mark run font hint of
# "eastAsia" - matching what Word produces when two Tab presses are
# typed into an empty East-Asian-aware paragraph.

$d = $word.ActiveDocument
$target = $d.Paragraphs.Last.Range

$openXmlPackage = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
          </w:p>
          <w:sectPr>
            <w:pgSz w:w="11906" w:h="16838"/>
          </w:sectPr>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($openXmlPackage)
